$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.794.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.71%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.291.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.29%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.04%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'299.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'97.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.65%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.512"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.85%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.07%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.505"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.78%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'35.82"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.95%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -0.55%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +0.74%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'17.77"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.27%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'6.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.52%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'2.648.49"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.25%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'2.288.80"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -2.86%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.775"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.53%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'42.722.85"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.70%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'12.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -5.45%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.0₃0904"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.65%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'6.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.43%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'67.77"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.74%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'241.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.69%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.31%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +0.05%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.01%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'4.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.52%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'25.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.51%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'165.88"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.80%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'2.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.77%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'9.02"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -2.01%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'32.78"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -3.59%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D34").Value = "'4.72"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -5.73%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'5.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -3.61%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'17.09"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -6.56%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'2.37"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.37%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.0684"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -2.08%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -2.17%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -4.63%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -1.68%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.110"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.24%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'2.015.54"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.96%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D45").Value = "'10.09"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.43%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -3.82%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'17.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.81%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.77"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -3.22%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'2.516.34"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.20%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'53.07"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -3.86%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'2.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -6.66%  "
$ws.Range("E51").Style = "Normal"
